# Negative Jump Detailed Cases.xlsx - apply the "TAI1 - RTAI / UTC1" table update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- existing table (rows 67-70): D70 "N / A" -> "Handled by C-A" ---
$ws.Range("D70").Value = "Handled by C-A"

# --- new mirrored table for "TAI1 - RTAI / UTC1" (rows 84-88) ---
$ws.Range("A84").Value = "TAI1 - RTAI / UTC1"
$ws.Range("B84").Value = "A - Ambiguous Region Before the Jump"
$ws.Range("C84").Value = "B - Ambiguous Region After the Jump"
$ws.Range("D84").Value = "C - Not in Ambiguous Region"

$ws.Range("A85").Value = "A - In Ambiguous Region"
$ws.Range("B85").Value = "Use pre Jump value of UTC2"
$ws.Range("C85").Value = "Use post Jump value of UTC2"
$ws.Range("D85").Value = "Use post Jump value of UTC2 (since UTC1 had to be after ambiguous regions)"

$ws.Range("A86").Value = "B - Not in Ambiguous Region"
$ws.Range("B86").Value = "No adjustment"
$ws.Range("C86").Value = "Make a positive adjustment to TAI1 - RTAI"
$ws.Range("D86").Value = "Trivial case"

$ws.Range("A87").Value = "C - In Other Ambiguous Region"
$ws.Range("B87").Value = "Use post jump value of UTC2"
$ws.Range("C87").Value = "Add .10 to TAI1 - RTAI and use post jump value of UTC2"
$ws.Range("D87").Value = "Handled by C-A"

# B88 is a blank spacer cell (matches the blank cell under B71 in the original
# table); stamp then clear it so a row/cell entry exists, and share the same
# "no-op alignment" formatting as B86/B87 (mirrors the B69/B70/B71 column above).
$ws.Range("B88").Value = "x"
$ws.Range("B86:B88").IndentLevel = 0
$ws.Range("B88").ClearContents() | Out-Null

# --- view state: active selection moved along with the new content ---
$ws.Range("D71").Select() | Out-Null
